# Fix query logging: add "Issue" column (F) and append newly logged rows 8-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell F1, matching the style of the existing header row ---
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F1").Value = "Issue"

# --- Newly logged helpline rows (A:E populated like the existing log rows) ---
$ws.Range("A8").Value = "2025-06-13T18:57:46.456198"
$ws.Range("B8").Value = "Abhinab Kumar"
$ws.Range("D8").Value = 9386776366
$ws.Range("E8").Value = "Name: Abhinab Kumar, Mobile Number: 9386776366"

$ws.Range("A9").Value = "2025-06-13T18:59:17.712704"
$ws.Range("B9").Value = "Abhinab"
$ws.Range("D9").Value = 8787878787
$ws.Range("E9").Value = "Name: Abhinab, Mobile Number: 8787878787"

$ws.Range("A10").Value = "2025-06-13T19:00:44.873436"
$ws.Range("B10").Value = "Abhinab"
$ws.Range("D10").Value = 9386776366
$ws.Range("E10").Value = "unable to upload docs"

$ws.Range("A11").Value = "2025-06-13T19:01:08.962484"
$ws.Range("B11").Value = "Akash Kr"
$ws.Range("D11").Value = 9386776366
$ws.Range("E11").Value = "Report an issue"

$ws.Range("A12").Value = "2025-06-13T19:02:56.754380"
$ws.Range("B12").Value = "Abhinab Kumare"
$ws.Range("D12").Value = 9386776366
$ws.Range("E12").Value = "Please describe the issue you are facing."

$ws.Range("A13").Value = "2025-06-13T19:03:12.631977"
$ws.Range("B13").Value = "Abhinab Kumare"
$ws.Range("D13").Value = 9386776366
$ws.Range("E13").Value = "Unable to upload docs"

$ws.Range("A14").Value = "2025-06-13T19:08:56.093580"
$ws.Range("B14").Value = "Abhinab Kumar"
$ws.Range("D14").Value = 9386776366
$ws.Range("E14").Value = "Please also describe your issue."

# --- Rows that now use the new "Issue" (F) column instead of "Query" (E) ---
$ws.Range("A15").Value = "2025-06-13T19:32:20.394277"
$ws.Range("B15").Value = "Nishant Kumar"
$ws.Range("C15").Value = "ABCZYX1"
$ws.Range("F15").Value = "Unable to get money"

$ws.Range("A16").Value = "2025-06-13T19:41:24.522473"
$ws.Range("B16").Value = "Aditya Thakur"
$ws.Range("D16").Value = 9386776366
$ws.Range("F16").Value = "Report an Issue"

$ws.Range("A17").Value = "2025-06-13T19:43:02.460239"
$ws.Range("B17").Value = "Aditya Thakue"
$ws.Range("D17").Value = 9386776366
$ws.Range("F17").Value = "unable to do anything"

$ws.Range("A18").Value = "2025-06-13T19:46:53.270003"
$ws.Range("B18").Value = "Bittu Kumar"
$ws.Range("C18").Value = "9892BX"
$ws.Range("F18").Value = "Unable to open docs"

$ws.Range("A19").Value = "2025-06-13T19:51:08.681914"
$ws.Range("B19").Value = "Abhinab Kumar"
$ws.Range("C19").Value = "87654BX"
$ws.Range("F19").Value = "Unable to log in"

Write-Output "applied helpline log update"
